$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-11 from 2023-09-01 (45170)
# to 2023-09-05 (45174), keeping the existing date formatting.
$ws.Range("C2:C11").Value = 45174
